$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-10 to new serial date value 45204
foreach ($row in 2..10) {
    $ws.Cells.Item($row, 3).Value = 45204
}
